$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Modelo" in F1, matching the style of the existing header row
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"

# Update existing metric values in row 2
$ws.Range("B2").Value = 0.2138012735004911
$ws.Range("C2").Value = 0.9842035830415448
$ws.Range("D2").Value = 0.3492862600957318

# Add new model name value in F2
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor(n_estimators=100))])"
